$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date-type/report-date text cells (keep as text)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Update the numeric financial figures for row 2
$ws.Range("O2").Value = 49051656.55
$ws.Range("P2").Value = 341591337.46
$ws.Range("Q2").Value = 299781032.66
$ws.Range("R2").Value = -27.4971726073
$ws.Range("S2").Value = 264608664.03
$ws.Range("T2").Value = 264608664.03
$ws.Range("U2").Value = -26.7914839042
$ws.Range("V2").Value = 7495410.22
$ws.Range("W2").Value = 12497280.4
$ws.Range("X2").Value = -849097.29
$ws.Range("Y2").Value = 56440136.35
$ws.Range("Z2").Value = 56530098.95
$ws.Range("AA2").Value = 7478442.4
$ws.Range("AG2").Value = 2989877.93
$ws.Range("AP2").Value = -29.2108396714
$ws.Range("AQ2").Value = -14.602149930832
$ws.Range("AR2").Value = -11.689896640096
$ws.Range("AS2").Value = 37050556.55
$ws.Range("AT2").Value = -9.527612356043999
